# Apply "Listing users, towns and categories" edit:
# - Mark Admin Edit Ad (row 37), Admin Delete Ad (row 38), Admin List Users (row 39),
#   Admin List Categories (row 43) and Admin List Towns (row 47) as "Yes" in column E.
# - Bump GitHub commit counts for the two students (C8, C9).
# - Update the active selection/scroll position of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GitHub commit counts.
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 18

# Mark these admin feature rows as implemented ("Yes").
$ws.Range("E37").Value = "Yes"
$ws.Range("E38").Value = "Yes"
$ws.Range("E39").Value = "Yes"
$ws.Range("E43").Value = "Yes"
$ws.Range("E47").Value = "Yes"

# Move the view/selection like the author did while editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("J11").Select()
